# Update balance-analysis results on the "GLOBAL RESULTS" sheet and the
# "LANDING GEARS" sheet with refreshed Xcg/Ycg/Zcg computations.

$wb = $excel.ActiveWorkbook

$wsGlobal = $wb.Worksheets.Item("GLOBAL RESULTS")
$wsGlobal.Range("C2").Value = 45.89293685237998
$wsGlobal.Range("C3").Value = 12.240997717128838
$wsGlobal.Range("C4").Value = 0.6576665053638682
$wsGlobal.Range("C6").Value = 26.36362870317379
$wsGlobal.Range("C7").Value = 11.791351849019662
$wsGlobal.Range("C8").Value = 0.721704678058392
$wsGlobal.Range("C10").Value = 26.36362870317379
$wsGlobal.Range("C11").Value = 11.791351849019662
$wsGlobal.Range("C12").Value = 0.721704678058392
$wsGlobal.Range("C14").Value = 30.376468346802522
$wsGlobal.Range("C15").Value = 11.883744101287732
$wsGlobal.Range("C16").Value = 0.46534746772689106
$wsGlobal.Range("C18").Value = 19.737044834211375
$wsGlobal.Range("C19").Value = 11.638780337853301
$wsGlobal.Range("C20").Value = 0.6736251083832776

$wsGears = $wb.Worksheets.Item("LANDING GEARS")
$wsGears.Range("C2").Value = 12.299024241711926
